$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new KP Training Acc / KP Test Acc values for rows 2-5
# Row 2: iris 77/73
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1

# Row 3: iris 113/37
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1

# Row 4: sonar 157/51
$ws.Range("B4").Value = 0.5414
$ws.Range("C4").Value = 0.5098

# Row 5: sonar 116/92
$ws.Range("B5").Value = 0.5517
$ws.Range("C5").Value = 0.5109

# Update the selection to A6, matching the saved view state
$ws.Range("A6").Select()
